# Insert a new data row at row 46 (pushing the existing rows 46-49 down to
# 47-50), then populate the new row with its data. This mirrors the diff:
# a weekly update that adds one more price record for "Arveja Verde".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 46:49 down to 47:50, freeing up row 46 for the new record.
$ws.Rows.Item(46).Insert()

# Populate the newly inserted row 46 with the new record's data.
$ws.Cells.Item(46, 1).Value = 11
$ws.Cells.Item(46, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(46, 3).Value = "Bíobío"
$ws.Cells.Item(46, 4).Value = 45218
$ws.Cells.Item(46, 4).NumberFormat = $ws.Cells.Item(47, 4).NumberFormat
$ws.Cells.Item(46, 5).Value = 8
$ws.Cells.Item(46, 6).Value = 100112022
$ws.Cells.Item(46, 7).Value = "Arveja Verde"
$ws.Cells.Item(46, 8).Value = "Perfection"
$ws.Cells.Item(46, 9).Value = "Primera"
$ws.Cells.Item(46, 10).Value = 50
$ws.Cells.Item(46, 11).Value = 28000
$ws.Cells.Item(46, 12).Value = 28000
$ws.Cells.Item(46, 13).Value = 28000
$ws.Cells.Item(46, 14).Value = "`$/malla 25 kilos"
$ws.Cells.Item(46, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(46, 16).Value = 1120
$ws.Cells.Item(46, 17).Value = 25
$ws.Cells.Item(46, 18).Value = "Hortaliza"
